$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the b.md row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 02:44:29"

# --- zh-cn sheet: row 3 is the b.md row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Copy/PasteSpecial from a cell that already holds "False" as a shared
# string so the target cell keeps the string type instead of becoming
# a native boolean cell.
$wsZhCn.Range("F2").Copy()
$wsZhCn.Range("F3").PasteSpecial(-4163)
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-02 02:44:24"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/506073de53d7920d244c72963c98ecbc46c5c8b5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc5f32b4e032fa7149e3e228b234c808a55b6715/e2e/b.md."
# ColumnWidth (character units) maps to OOXML "width" with a fixed offset
# in this runtime; 39.14 yields the target raw width of 40.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.14

# --- de-de sheet: row 3 is the b.md row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F2").Copy()
$wsDeDe.Range("F3").PasteSpecial(-4163)
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-02 02:44:29"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/506073de53d7920d244c72963c98ecbc46c5c8b5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc5f32b4e032fa7149e3e228b234c808a55b6715/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.14
